# Updated cryptos list on Sat Jun 15 15:31:52 UTC 2024 with GitHub Actions
# Refreshes Price (D) / Volume(1h) (E) figures for each coin row, including
# two coin-pair swaps (rows 29<->30 and 34<->35) and one coin replacement
# (row 51: EnergySwap -> SuiNetwork). Numeric-looking text values are
# written with a leading "'" so Excel keeps them as text instead of
# auto-converting to numbers (matching the source data's inline-string
# formatting).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.351.37'
$ws.Range("E2").Value = '  -0.89%  '

$ws.Range("D3").Value = '3.578.82'
$ws.Range("E3").Value = '  +2.58%  '

$ws.Range("D4").Value = "'" + '1.00'
$ws.Range("E4").Value = '  +0.17%  '

$ws.Range("D5").Value = "'" + '608.23'
$ws.Range("E5").Value = '  +0.43%  '

$ws.Range("D6").Value = "'" + '145.09'
$ws.Range("E6").Value = '  -0.13%  '

$ws.Range("D7").Value = '3.575.95'
$ws.Range("E7").Value = '  +2.54%  '

$ws.Range("E8").Value = '  +0.14%  '

$ws.Range("D9").Value = "'" + '0.484'
$ws.Range("E9").Value = '  +1.51%  '

$ws.Range("E10").Value = '  -2.53%  '

$ws.Range("D11").Value = "'" + '8.01'
$ws.Range("E11").Value = '  +0.81%  '

$ws.Range("D12").Value = "'" + '0.412'
$ws.Range("E12").Value = '  -1.19%  '

$ws.Range("D13").Value = '4.196.78'
$ws.Range("E13").Value = '  +3.11%  '

$ws.Range("D14").Value = "'" + '0.0000209'
$ws.Range("E14").Value = '  -1.67%  '

$ws.Range("D15").Value = "'" + '30.33'
$ws.Range("E15").Value = '  -2.44%  '

$ws.Range("D16").Value = '3.585.42'
$ws.Range("E16").Value = '  +3.16%  '

$ws.Range("D17").Value = '66.454.34'
$ws.Range("E17").Value = '  -0.74%  '

$ws.Range("D18").Value = "'" + '11.69'
$ws.Range("E18").Value = '  +10.18%  '

$ws.Range("E19").Value = '  -1.28%  '

$ws.Range("D20").Value = "'" + '6.23'
$ws.Range("E20").Value = '  -0.74%  '

$ws.Range("D21").Value = "'" + '15.00'
$ws.Range("E21").Value = '  -2.21%  '

$ws.Range("D22").Value = "'" + '431.17'
$ws.Range("E22").Value = '  +0.56%  '

$ws.Range("D23").Value = "'" + '0.609'
$ws.Range("E23").Value = '  +1.31%  '

$ws.Range("D24").Value = "'" + '78.76'
$ws.Range("E24").Value = '  -0.68%  '

$ws.Range("D25").Value = '3.729.03'

$ws.Range("E26").Value = '  +0.05%  '

$ws.Range("E27").Value = '  +4.35%  '

$ws.Range("D28").Value = "'" + '8.09'
$ws.Range("E28").Value = '  -0.26%  '

$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").Value = "'" + '2.53'
$ws.Range("E29").Value = '  +1.73%  '

$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").Value = "'" + '9.28'
$ws.Range("E30").Value = '  -4.31%  '

$ws.Range("E31").Value = '  -0.16%  '

$ws.Range("E32").Value = '  -3.33%  '

$ws.Range("E33").Value = '  -3.21%  '

$ws.Range("B34").Value = 'RenzoRestakedETH'
$ws.Range("C34").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D34").Value = '3.582.56'
$ws.Range("E34").Value = '  +2.90%  '

$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D35").Value = "'" + '25.50'
$ws.Range("E35").Value = '  +0.70%  '

$ws.Range("D36").Value = "'" + '1.76'
$ws.Range("E36").Value = '  -0.42%  '

$ws.Range("E37").Value = '  -0.04%  '

$ws.Range("E38").Value = '  -0.26%  '

$ws.Range("D39").Value = "'" + '5.67'
$ws.Range("E39").Value = '  -0.59%  '

$ws.Range("D40").Value = "'" + '1.00'
$ws.Range("E40").Value = '  +0.22%  '

$ws.Range("D41").Value = "'" + '172.05'
$ws.Range("E41").Value = '  -1.70%  '

$ws.Range("D42").Value = "'" + '0.0860'
$ws.Range("E42").Value = '  -3.38%  '

$ws.Range("D43").Value = "'" + '5.29'
$ws.Range("E43").Value = '  -0.36%  '

$ws.Range("D44").Value = "'" + '0.898'
$ws.Range("E44").Value = '  +0.79%  '

$ws.Range("E45").Value = '  -2.38%  '

$ws.Range("D46").Value = "'" + '45.92'
$ws.Range("E46").Value = '  -0.66%  '

$ws.Range("E47").Value = '  +1.76%  '

$ws.Range("D48").Value = "'" + '26.04'
$ws.Range("E48").Value = '  -5.79%  '

$ws.Range("D49").Value = "'" + '2.41'
$ws.Range("E49").Value = '  +1.05%  '

$ws.Range("D50").Value = "'" + '7.16'
$ws.Range("E50").Value = '  -1.96%  '

$ws.Range("B51").Value = 'SuiNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D51").Value = "'" + '0.954'
$ws.Range("E51").Value = '  -2.12%  '
